# ENG-110 "Intro.docx" edit
#
# 1) In the bullet list, split the sentence "Over these past few years,
#    there's been a lot of protest against the censorship." so that
#    "protest against" and "the censorship" each get their own run,
#    bracketed by grammar-check proofErr (gramStart/gramEnd) markers.
# 2) In the merged recap paragraph (page 2), reword
#    "Diverse literature is extremely important for developing empathy
#    in students. Developmental malleability is crucial..." into
#    "Diverse literature has proven important for developing empathy in
#    students and developmental malleability is crucial..." (broken into
#    several runs) and apply the same protest-against/the-censorship
#    run split there too.

$d = $word.ActiveDocument

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgClose = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

$rPr = '<w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr>'

function Escape-Xml([string]$text) {
    $t = $text -replace '&', '&amp;'
    $t = $t -replace '<', '&lt;'
    $t = $t -replace '>', '&gt;'
    return $t
}

function New-Run([string]$text) {
    $escaped = Escape-Xml $text
    if ($text.Length -eq 0 -or $text[0] -eq ' ' -or $text[-1] -eq ' ') {
        return '<w:r>' + $rPr + '<w:t xml:space="preserve">' + $escaped + '</w:t></w:r>'
    }
    return '<w:r>' + $rPr + '<w:t>' + $escaped + '</w:t></w:r>'
}

function New-GramRun([string]$text) {
    return '<w:proofErr w:type="gramStart"/>' + (New-Run $text) + '<w:proofErr w:type="gramEnd"/>'
}

# The split used both times "... a lot of protest against the censorship..."
# turns into "...a lot of [protest against] [ ] [the censorship]..."
function New-ProtestSplit() {
    return (New-GramRun "protest against") + (New-Run " ") + (New-GramRun "the censorship")
}

# ---------------------------------------------------------------------
# Change 1: bullet paragraph.
# ---------------------------------------------------------------------
$oldSentence = "Over these past few years, there" + [char]0x2019 + "s been a lot of protest against the censorship."

$bulletPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq $oldSentence) {
        $bulletPara = $p
        break
    }
}

$innerXml1 = (New-Run ("Over these past few years, there" + [char]0x2019 + "s been a lot of ")) + `
    (New-ProtestSplit) + (New-Run ".")

$r1 = $d.Range($bulletPara.Range.Start, $bulletPara.Range.End - 1)
$r1.InsertXML($pkgOpen + '<w:p>' + $innerXml1 + '</w:p>' + $pkgClose)

# ---------------------------------------------------------------------
# Change 2: merged recap paragraph on page 2. We rebuild the whole
# paragraph (InsertXML only splices correctly when given the complete
# paragraph range) but keep the two runs that bracket the sentence
# ("...case." and "By banning...") byte-for-byte identical to the
# original, including their rsid/page-break markup.
# ---------------------------------------------------------------------
$mergedPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains("Diverse literature is extremely important")) {
        $mergedPara = $p
        break
    }
}

$full = $mergedPara.Range
$txt = $full.Text

$markerA = "Stop W.O.K.E."
$idxAEnd = $txt.IndexOf($markerA) + $markerA.Length

$oldMerged = " Diverse literature is extremely important for developing empathy in students. Developmental malleability is crucial for efficient learning. Over these past few years, there" + [char]0x2019 + "s been a lot of protest against the censorship. Many school board meetings have been held whether books on LGBTQ+ topics should be kept in the classroom. "
$idxMergedStart = $txt.IndexOf($oldMerged)
$idxMergedEnd = $idxMergedStart + $oldMerged.Length

$runAText = $txt.Substring(0, $idxAEnd)
$runA2Text = $txt.Substring($idxAEnd, $idxMergedStart - $idxAEnd)
$runCText = $txt.Substring($idxMergedEnd)

# Run A keeps its original <w:lastRenderedPageBreak/> before the text.
$runA = '<w:r>' + $rPr + '<w:lastRenderedPageBreak/><w:t>' + (Escape-Xml $runAText) + '</w:t></w:r>'
# Run A2 and run C originally carry w:rsidR="001C0DB3".
$runA2 = '<w:r w:rsidR="001C0DB3">' + $rPr + '<w:t xml:space="preserve">' + (Escape-Xml $runA2Text) + '</w:t></w:r>'
$runC = '<w:r w:rsidR="001C0DB3">' + $rPr + '<w:t>' + (Escape-Xml $runCText) + '</w:t></w:r>'

$innerXml2 = (New-Run " Diverse literature ") + (New-Run "has proven important") + `
    (New-Run " for developing empathy in students") + (New-Run " and") + (New-Run " ") + `
    (New-Run "d") + `
    (New-Run ("evelopmental malleability is crucial for efficient learning. Over these past few years, there" + [char]0x2019 + "s been a lot of ")) + `
    (New-ProtestSplit) + `
    (New-Run ". Many school board meetings have been held whether books on LGBTQ+ topics should be kept in the classroom. ")

$allRuns = $runA + $runA2 + $innerXml2 + $runC

$r2 = $d.Range($full.Start, $full.End - 1)
$r2.InsertXML($pkgOpen + '<w:p>' + $allRuns + '</w:p>' + $pkgClose)
